# Update the cryptos price/volume table with the latest scraped values.
# Values are prefixed with a leading apostrophe so Excel stores them as text
# (matching the original inlineStr cells) instead of auto-converting
# numeric-looking strings (e.g. "206.50", "0.253") into floating point
# numbers, which would drop trailing zeros / introduce FP rounding noise.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.662.25"
$ws.Range("E2").Value = "'  -0.89%  "
$ws.Range("E3").Value = "'  -3.30%  "
$ws.Range("E4").Value = "'  +0.28%  "
$ws.Range("D5").Value = "'206.50"
$ws.Range("E5").Value = "'  -2.51%  "
$ws.Range("E6").Value = "'  -3.15%  "
$ws.Range("E7").Value = "'  +0.30%  "
$ws.Range("D8").Value = "'22.31"
$ws.Range("E8").Value = "'  -4.77%  "
$ws.Range("D9").Value = "'0.253"
$ws.Range("E9").Value = "'  -1.45%  "
$ws.Range("E10").Value = "'  -3.41%  "
$ws.Range("E11").Value = "'  -1.69%  "
$ws.Range("D12").Value = "'1.807.18"
$ws.Range("E12").Value = "'  -3.31%  "
$ws.Range("D13").Value = "'1.585.51"
$ws.Range("E13").Value = "'  -3.19%  "
$ws.Range("E14").Value = "'  -4.09%  "
$ws.Range("D15").Value = "'0.530"
$ws.Range("E15").Value = "'  -5.70%  "
$ws.Range("D16").Value = "'27.633.88"
$ws.Range("E16").Value = "'  -0.97%  "
$ws.Range("D17").Value = "'63.03"
$ws.Range("E17").Value = "'  -3.58%  "
$ws.Range("D18").Value = "'217.84"
$ws.Range("E18").Value = "'  -4.89%  "
$ws.Range("E19").Value = "'  -3.63%  "
$ws.Range("D20").Value = "'7.31"
$ws.Range("E20").Value = "'  -5.04%  "
$ws.Range("E21").Value = "'  +0.26%  "
$ws.Range("E22").Value = "'  -4.82%  "
$ws.Range("D23").Value = "'9.51"
$ws.Range("E23").Value = "'  -5.39%  "
$ws.Range("E24").Value = "'  -5.00%  "
$ws.Range("D25").Value = "'153.59"
$ws.Range("E25").Value = "'  -1.16%  "
$ws.Range("E26").Value = "'  +0.30%  "
$ws.Range("D27").Value = "'6.69"
$ws.Range("E27").Value = "'  -2.80%  "
$ws.Range("D28").Value = "'15.06"
$ws.Range("E28").Value = "'  -3.11%  "
$ws.Range("E29").Value = "'  -4.36%  "
$ws.Range("E30").Value = "'  -2.84%  "
$ws.Range("E31").Value = "'  -3.60%  "
$ws.Range("E32").Value = "'  -5.41%  "
$ws.Range("D33").Value = "'1.375.07"
$ws.Range("E33").Value = "'  -1.57%  "
$ws.Range("D34").Value = "'2.94"
$ws.Range("E34").Value = "'  -5.42%  "
$ws.Range("D35").Value = "'1.51"
$ws.Range("E35").Value = "'  -5.68%  "
$ws.Range("D36").Value = "'0.967"
$ws.Range("E36").Value = "'  -4.79%  "
$ws.Range("E37").Value = "'  -1.32%  "
$ws.Range("E38").Value = "'  -3.15%  "
$ws.Range("E39").Value = "'  -3.74%  "
$ws.Range("D40").Value = "'0.817"
$ws.Range("E40").Value = "'  -3.96%  "
$ws.Range("E41").Value = "'  +0.32%  "
$ws.Range("D42").Value = "'0.980"
$ws.Range("E42").Value = "'  -2.69%  "
$ws.Range("E43").Value = "'  -3.29%  "
$ws.Range("B44").Value = "'Aave"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'63.70"
$ws.Range("E44").Value = "'  -3.48%  "
$ws.Range("B45").Value = "'MXToken"
$ws.Range("C45").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.17"
$ws.Range("E45").Value = "'  +1.33%  "
$ws.Range("D46").Value = "'5.20"
$ws.Range("E46").Value = "'  -4.54%  "
$ws.Range("D47").Value = "'1.719.04"
$ws.Range("E47").Value = "'  -3.22%  "
$ws.Range("D48").Value = "'87.62"
$ws.Range("E48").Value = "'  -1.21%  "
$ws.Range("E49").Value = "'  -2.52%  "
$ws.Range("D50").Value = "'0.0972"
$ws.Range("E50").Value = "'  -5.54%  "
$ws.Range("E51").Value = "'  -1.50%  "
